# para agregar codigo de barras
# Updates patient identification data on the hospital intake form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name block (row 6) ---
$ws.Range("A6").Value = "García"
$ws.Range("C6").Value = "Figueroa"
$ws.Range("E6").Value = "kevin"
$ws.Range("G6").Value = "Estuardo"
$ws.Range("I6").Value = "/201757568"

# --- Dirección actual (row 8 label values / row 10 data) ---
$ws.Range("H8").Value = "Gutemala"
$ws.Range("J8").Value = "47000538"

$ws.Range("A10").Value = "9a. ave 4-71 "
$ws.Range("D10").Value = "Zona 19 colonia la florida"
$ws.Range("F10").Value = "guatemala"
$ws.Range("H10").Value = "guatemala"
$ws.Range("J10").Value = "47000538"

# --- Fecha de nacimiento / edad / lugar / sexo (row 12) ---
$ws.Range("A12").Value = "23/11/1992"
$ws.Range("F12").Value = " 24"
$ws.Range("H12").Value = "guatemala"
$ws.Range("J12").Value = "MASCULINO"

# --- Estado civil / ocupación / nacionalidad / cédula (row 14) ---
$ws.Range("A14").Value = "Soltero"
$ws.Range("D14").Value = " estudiante "
$ws.Range("F14").Value = "guatemalteco"
$ws.Range("H14").Value = "2424583430101"

# --- Nombre del Padre / Madre (row 18) ---
$ws.Range("A18").Value = "Manolo García"
$ws.Range("F18").Value = "Aura Leticia Figueroa"

# --- Emergency contact (row 20) ---
$ws.Range("A20").Value = "Aura Leticia Figueroa"
$ws.Range("F20").Value = "madre "
$ws.Range("H20").Value = ""
$ws.Range("J20").Value = "47000538"

# --- Fecha / hora de ingreso (row 24) ---
$ws.Range("A24").Value = "17/11/2017"
$ws.Range("C24").Value = "11:59:50"
